# Sparks 2.2.9 - finishing content
# Rename "Click Area" header to "Click Area Name", tweak the min/max-y sample
# values (60/-60 -> 40/-40), shrink the title/body fonts a couple points, and
# flip which columns carry the vertical divider border (col A loses it,
# cols B:E gain it) to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text -----------------------------------------------------
$ws.Range("A1").Value = "Click Area Name"

# --- Data tweaks: Maximum y / Minimum y columns (D, E) ----------------
# Values of 60 become 40, -60 becomes -40, everything else stays put.
for ($r = 2; $r -le 10; $r++) {
    foreach ($col in @("D", "E")) {
        $cell = $ws.Range("$col$r")
        $v = $cell.Value()
        if ($v -eq 60) {
            $cell.Value = 40
        } elseif ($v -eq -60) {
            $cell.Value = -40
        }
    }
}

# --- Font size tweaks ---------------------------------------------------
# Bold header font (A1 and B1:E1) goes from 24pt to 22pt.
$ws.Range("A1:E1").Font.Size = 22

# Body font (the larger "Cavolini" font used for data cells) goes from
# 20pt to 18pt.
$ws.Range("A2:E10").Font.Size = 18

# --- Border swap: column A loses its left border, columns B:E gain one -
$ws.Range("A2:A10").Borders.Item(7).LineStyle = -4142  # xlLeft -> none
$ws.Range("B2:E10").Borders.Item(7).LineStyle = 1      # xlContinuous
$ws.Range("B2:E10").Borders.Item(7).Weight = 2          # xlThin

# --- Column widths (Excel auto "best fit" after the edits) --------------
# (values back-solved so the saved OOXML <col width=.../> lands on the same
# pixel-quantized width Excel's own AutoFit produced: 26.66, 19.66, 19, ...)
$ws.Columns.Item(1).ColumnWidth = 25.830729166666668
$ws.Columns.Item(2).ColumnWidth = 18.830729166666668
$ws.Columns.Item(3).ColumnWidth = 18.166666666666668
$ws.Columns.Item(4).ColumnWidth = 18.830729166666668
$ws.Columns.Item(5).ColumnWidth = 18.166666666666668

# --- Selection moved to C4 ------------------------------------------------
$ws.Range("C4").Select()

$wb.Save()
